# Replace the two-digit multiplication problems/answers with the new values.
$d = $word.ActiveDocument

$replacements = @(
    @{old = "97×13=1261"; new = "70×89=6230"},
    @{old = "45×95=4275"; new = "29×40=1160"},
    @{old = "19×64=1216"; new = "97×65=6305"},
    @{old = "81×56=4536"; new = "84×13=1092"},
    @{old = "59×48=2832"; new = "98×51=4998"},
    @{old = "25×72=1800"; new = "78×49=3822"},
    @{old = "53×59=3127"; new = "69×21=1449"},
    @{old = "75×26=1950"; new = "75×62=4650"},
    @{old = "42×12=504";  new = "90×80=7200"},
    @{old = "98×35=3430"; new = "57×29=1653"},
    @{old = "50×17=850";  new = "45×70=3150"},
    @{old = "84×43=3612"; new = "35×66=2310"},
    @{old = "69×13=897";  new = "32×59=1888"},
    @{old = "24×18=432";  new = "93×32=2976"},
    @{old = "20×92=1840"; new = "25×17=425"},
    @{old = "61×38=2318"; new = "34×60=2040"},
    @{old = "82×79=6478"; new = "70×94=6580"},
    @{old = "78×43=3354"; new = "12×45=540"},
    @{old = "78×34=2652"; new = "93×95=8835"},
    @{old = "26×91=2366"; new = "47×45=2115"},
    @{old = "85×54=4590"; new = "57×76=4332"},
    @{old = "31×27=837";  new = "61×85=5185"},
    @{old = "52×42=2184"; new = "29×13=377"},
    @{old = "98×67=6566"; new = "14×40=560"},
    @{old = "15×47=705";  new = "95×53=5035"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}

$d.Save()
